$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlotSettings")

# Remove the "combined transparency" row entirely (was row 4); every row
# below it shifts up by one.
$ws.Rows(4).Delete()

# Rename the remaining setting labels in column A (new row numbers, after
# the deletion above) to their shorter forms.
$ws.Range("A2").Value = "Linewidth"
$ws.Range("A3").Value = "Arrowsize"
$ws.Range("A4").Value = "FontSize"
$ws.Range("A5").Value = "NodeSize"
$ws.Range("A6").Value = "NodeColor"
